# Update handback status timestamps (regenerate report for handback)
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 05:03:56"
$wsZhCn.Range("H2").Value = "2016-03-13 05:04:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 05:04:03"
$wsDeDe.Range("H2").Value = "2016-03-13 05:04:24"
